$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Idle time (column R) = 1 - rho, for data rows 2 through 18
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 18)  # column R = 18
    $old = $cell.Value()
    $cell.Value = 1 - $old
}

# Update the active selection as reflected in the saved file
$ws.Range("U14").Select()
